$p = $ppt.ActivePresentation

$oldDate = "9/15/2025"
$newDate = "9/16/2025"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDate = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDate = $true
                }
            } catch {
                $isDate = $false
            }
            if (-not $isDate -and $shp.Name -like "Date Placeholder*") {
                $isDate = $true
            }
            if ($isDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master date placeholder
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout (CustomLayout) date placeholder
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master date placeholder
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
